$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows for Bradley Beal / Tari Eason / De'Andre Hunter / Dyson Daniels
# were reordered: Bradley Beal (row 13) moved down below Dyson Daniels
# (to row 16), shifting the other three rows up by one.
# Net effect on cell contents, row by row:

$ws.Range("A13").Value = "Tari Eason"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Houston Rockets"

$ws.Range("A14").Value = "De'Andre Hunter"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "Atlanta Hawks"

$ws.Range("A15").Value = "Dyson Daniels"
$ws.Range("B15").Value = "PG,SG,SF"
$ws.Range("C15").Value = "Atlanta Hawks"

$ws.Range("A16").Value = "Bradley Beal"
$ws.Range("B16").Value = "PG,SG,SF"
$ws.Range("C16").Value = "Phoenix Suns"
